# Auto-generated: applies cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with simple in-place Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "67.573.70"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").Value = "3.776.10"
$ws.Range("E3").Value = "  +8.21%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'419.11"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'132.63"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "3.755.11"
$ws.Range("E7").Value = "  +7.86%  "
$ws.Range("D8").Value = "'0.646"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'0.767"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'0.183"
$ws.Range("E11").Value = "  +12.36%  "
$ws.Range("D12").Value = "'0.0000403"
$ws.Range("E12").Value = "  +51.69%  "
$ws.Range("D13").Value = "'42.51"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "'10.39"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").Value = "4.375.19"
$ws.Range("E15").Value = "  +8.21%  "
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "3.736.72"
$ws.Range("E17").Value = "  +6.57%  "
$ws.Range("D18").Value = "'20.44"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'13.21"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "'1.12"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "67.674.62"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("D22").Value = "'441.90"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "'15.75"
$ws.Range("E23").Value = "  +19.79%  "
$ws.Range("D24").Value = "'90.09"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "'3.07"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").Value = "'38.09"
$ws.Range("E26").Value = "  +12.44%  "
$ws.Range("D27").Value = "'3.32"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("D29").Value = "'5.09"
$ws.Range("E29").Value = "  +5.29%  "
$ws.Range("E30").Value = "  +6.09%  "
$ws.Range("D31").Value = "'12.54"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "'2.71"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "'7.12"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'41.21"
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("D36").Value = "'57.87"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'0.0488"
$ws.Range("E38").Value = "  -2.99%  "
$ws.Range("D44").Value = "'3.38"
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("D45").Value = "'148.61"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "'3.17"
$ws.Range("E46").Value = "  +24.47%  "
$ws.Range("D47").Value = "'2.10"
$ws.Range("E47").Value = "  +5.82%  "
$ws.Range("D48").Value = "'2.88"
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("D49").Value = "'2.61"
$ws.Range("E49").Value = "  -6.19%  "
$ws.Range("D50").Value = "'4.29"
$ws.Range("E50").Value = "  -4.56%  "
$ws.Range("D51").Value = "'0.304"
$ws.Range("E51").Value = "  -2.06%  "

# --- Rows 39-43: re-ranked coins (Coin/Link/Price/Volume all change; rank index in column A is unchanged) ---
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = "  +28.04%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.147"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0699"
$ws.Range("E41").Value = "  -3.30%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'27.61"
$ws.Range("E43").Value = "  +28.57%  "

